$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# C5 ("creation of nodejs web api") moves from the "in-progress" (theme7)
# fill to the same "done" (theme9) fill used by C4/C6/C7 - copy the format
# from C4 so the existing style is reused instead of creating a new one.
$ws1.Range("C4").Copy() | Out-Null
$ws1.Range("C5").PasteSpecial(-4122) | Out-Null

# Add the new weekly status worksheet right after Sheet1.
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$newSheet.Name = "Week 1 - Sprint 1"

$newSheet.Range("A1").Value = "w/c 30/01/2017"
$newSheet.Range("A1").Font.Bold = $true

$newSheet.Range("A3").Value = "tasks"
$newSheet.Range("A3").Font.Bold = $true

$newSheet.Range("B3").Value = 42765
$newSheet.Range("B3").NumberFormat = "mm-dd-yy"

$newSheet.Range("C3").Value = "tutor meeting: advice on proposal, notes in proposal document"

$newSheet.Range("A4").Formula = "=VALUE(Sheet1!A5)"
$newSheet.Range("C4").Value = "meeting with Tim B: need to talk to Steve about using Business Framework and coming into Sabisu, reviewed estimates, agreed to meet in either one or two weeks"

$newSheet.Range("C6").Value = "to do: tidy up of user interface/look more professional. Tidy up code, remove anything that is redundant. Create test plan"
$newSheet.Range("C5").Value = "work done: added delete and update functionality to web application, fixed bugs regarding click events, delete dialog added and overlay. Created stored procedures"

$newSheet.Columns.Item(1).ColumnWidth = 30.21875
$newSheet.Columns.Item(2).ColumnWidth = 10.5546875

# Match the saved selections from the source edit.
$ws1.Range("E8").Select() | Out-Null
$newSheet.Range("C5").Select() | Out-Null
